$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet's "key" cells (A8:A11) use a distinct font (Monaco, 11pt)
# compared to the default body font. Re-use that same formatting for the
# two rows being added below, matching the existing PUNISHMENT-* key rows.
$keyFontName = $ws.Range("A11").Font.Name
$keyFontSize = $ws.Range("A11").Font.Size

# Row 12: a blank "spacer" row (no value), but A12 still carries the
# key-cell formatting.
$ws.Range("A12").Font.Name = $keyFontName
$ws.Range("A12").Font.Size = $keyFontSize

# Row 14: new key row "PUNISHMENT-USERPARDON", same key-cell formatting as
# the other PUNISHMENT-* rows.
$ws.Range("A14").Font.Name = $keyFontName
$ws.Range("A14").Font.Size = $keyFontSize
$ws.Range("A14").Value = "PUNISHMENT-USERPARDON"

# Widen column A to fit the new, longer key and make it an explicit
# (no-longer best-fit) custom width.
$ws.Columns.Item(1).ColumnWidth = 35.333333

# Leave the active selection where the author's editing session ended.
$ws.Range("A19").Select() | Out-Null
